# Apio (Apio/Vega Central Mapocho) weekly price-list update.
#
# A new reporting date (2023-03-10, serial 44995) is inserted as a new pair
# of rows ("Primera"/"Segunda" quality) right above the existing 2022-08-25
# (serial 44798) pair, which in turn pushes every subsequent row down by two
# positions. The "Segunda" volume for the new date differs from the row it
# was cloned from (52 instead of 34). Everything else keeps the same
# relative order, just shifted down two rows, and the sheet's used range
# grows from A1:R406 to A1:R408.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the current 336:337 pair (Primera/Segunda @ 2022-08-25) and push
# everything from row 336 downward two rows lower — this is the COM
# equivalent of selecting rows 336:337, Copy, then Insert Copied Cells.
$ws.Range("A336:R337").Copy()
$ws.Range("A336:R337").Insert()

# Retarget the freshly-inserted pair to the new reporting date and fix up
# the "Segunda" volume for that date.
$ws.Range("D336").Value = 44995
$ws.Range("D337").Value = 44995
$ws.Range("J337").Value = 52
